$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9

# Row 5
$ws.Range("S5").Value = 1.63

# Row 6
$ws.Range("S6").Value = 1.47

# Row 7
$ws.Range("S7").Value = 1.54

# Row 8
$ws.Range("S8").Value = 1.47

# Row 9
$ws.Range("G9").Value = 1.85
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 1.36
$ws.Range("T9").Value = 3
$ws.Range("Z9").Value = 15
$ws.Range("AC9").Value = 11
$ws.Range("AE9").Value = 13
$ws.Range("AF9").Value = 41
$ws.Range("AS9").Value = 126
$ws.Range("AT9").Value = 3

# Row 10
$ws.Range("O10").Value = 1.17
$ws.Range("P10").Value = 5
$ws.Range("Q10").Value = 1.57
$ws.Range("R10").Value = 2.38

# Row 11
$ws.Range("G11").Value = 2.6
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 2.7
$ws.Range("J11").Value = 3.25
$ws.Range("K11").Value = 2.1
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 1.91
$ws.Range("X11").Value = 12
$ws.Range("AA11").Value = 21
$ws.Range("AC11").Value = 9.5
$ws.Range("AD11").Value = 6.5
$ws.Range("AG11").Value = 251
$ws.Range("AH11").Value = 8.5
$ws.Range("AL11").Value = 21
$ws.Range("AP11").Value = 23
$ws.Range("AR11").Value = 67
$ws.Range("AW11").Value = 4.75
$ws.Range("BA11").Value = 67
$ws.Range("BD11").Value = 126

# Row 15
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10

# Row 19
$ws.Range("G19").Value = 4.1
$ws.Range("H19").Value = 3.7
$ws.Range("I19").Value = 1.8
$ws.Range("L19").Value = 2.4
$ws.Range("O19").Value = 1.2
$ws.Range("P19").Value = 4.33
$ws.Range("Q19").Value = 1.7
$ws.Range("R19").Value = 2.1
$ws.Range("W19").Value = 15
$ws.Range("X19").Value = 23
$ws.Range("AK19").Value = 15
$ws.Range("AL19").Value = 13
$ws.Range("AM19").Value = 21
$ws.Range("AX19").Value = 9.5
$ws.Range("AY19").Value = 17
$ws.Range("AZ19").Value = 29
$ws.Range("BA19").Value = 41
$ws.Range("BB19").Value = 101

# Row 44
$ws.Range("G44").Value = 1.73
$ws.Range("H44").Value = 3.9
$ws.Range("I44").Value = 4.33
$ws.Range("J44").Value = 2.2
$ws.Range("N44").Value = 19
$ws.Range("O44").Value = 1.13
$ws.Range("P44").Value = 6
$ws.Range("Q44").Value = 1.44
$ws.Range("R44").Value = 2.7
$ws.Range("U44").Value = 1.44
$ws.Range("V44").Value = 2.63
$ws.Range("Z44").Value = 15
$ws.Range("AN44").Value = 4.33
$ws.Range("AO44").Value = 8.5
$ws.Range("AV44").Value = 41

# Row 45
$ws.Range("G45").Value = 1.27
$ws.Range("H45").Value = 5.75
$ws.Range("I45").Value = 9.5
$ws.Range("J45").Value = 1.67
$ws.Range("L45").Value = 8
$ws.Range("N45").Value = 21
$ws.Range("U45").Value = 1.83
$ws.Range("V45").Value = 1.83
$ws.Range("Z45").Value = 8.5
$ws.Range("AA45").Value = 11
$ws.Range("AE45").Value = 21
$ws.Range("AG45").Value = 251
$ws.Range("AI45").Value = 51
$ws.Range("AJ45").Value = 26
$ws.Range("AQ45").Value = 13
$ws.Range("AS45").Value = 101
$ws.Range("AW45").Value = 10
$ws.Range("AY45").Value = 41
$ws.Range("BA45").Value = 151
$ws.Range("BB45").Value = 251

# Row 48
$ws.Range("Q48").Value = 2.08
$ws.Range("R48").Value = 1.73

# Row 50
$ws.Range("BD50").Value = 151
